# Refresh crypto price/volume snapshot (GitHub Actions scheduled update).
# Source data is plain text (prices use "." as a thousands separator in some
# rows, so values like "71.029.59" are never valid numbers); a handful of
# updated prices DO look like plain numbers (e.g. "1.00", "17.54"). Writing
# those through .Value directly would make Excel auto-convert them to the
# numeric value 1 / 17.54 (dropping the literal text formatting), so for those
# cells we write them with a leading apostrophe to force text, then reset the
# cell style to Normal to drop the quote-prefix formatting Excel applies -
# leaving a plain text cell identical to the rest of the column.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "71.029.59"
$ws.Range("E2").Value = "  +0.17%  "
$ws.Range("D3").Value = "3.796.69"
$ws.Range("E3").Value = "  -0.79%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'699.77"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.34%  "
$ws.Range("D6").Value = "'169.92"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.16%  "
$ws.Range("D7").Value = "3.796.34"
$ws.Range("E7").Value = "  -0.76%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").Value = "'0.522"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.54%  "
$ws.Range("E10").Value = "  -1.37%  "
$ws.Range("E11").Value = "  +2.91%  "
$ws.Range("E12").Value = "  +4.23%  "
$ws.Range("D13").Value = "'0.0000249"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.90%  "
$ws.Range("D14").Value = "'36.15"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.45%  "
$ws.Range("D15").Value = "4.439.64"
$ws.Range("E15").Value = "  -0.78%  "
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "3.854.68"
$ws.Range("E16").Value = "  -0.93%  "
$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "71.274.57"
$ws.Range("E17").Value = "  +0.48%  "
$ws.Range("D18").Value = "'17.54"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Value = "'7.18"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.42%  "
$ws.Range("E20").Value = "  +0.00%  "
$ws.Range("D21").Value = "'515.33"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.85%  "
$ws.Range("D22").Value = "'10.40"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.38%  "
$ws.Range("D23").Value = "'0.713"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.59%  "
$ws.Range("D24").Value = "'83.57"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.88%  "
$ws.Range("D25").Value = "'0.0000140"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.81%  "
$ws.Range("D26").Value = "'12.54"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.29%  "
$ws.Range("D27").Value = "3.946.27"
$ws.Range("E27").Value = "  -0.92%  "
$ws.Range("D28").Value = "'10.21"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.63%  "
$ws.Range("E29").Value = "  +0.01%  "
$ws.Range("D30").Value = "'1.98"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.78%  "
$ws.Range("E31").Value = "  -3.41%  "
$ws.Range("D32").Value = "'2.25"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.00%  "
$ws.Range("D33").Value = "'7.27"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.32%  "
$ws.Range("D34").Value = "'29.12"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.49%  "
$ws.Range("E35").Value = "  -3.89%  "
$ws.Range("D36").Value = "'9.29"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.16%  "
$ws.Range("E37").Value = "  +1.32%  "
$ws.Range("D38").Value = "3.760.16"
$ws.Range("E38").Value = "  -0.86%  "
$ws.Range("D39").Value = "'6.56"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +10.03%  "
$ws.Range("E40").Value = "  -2.15%  "
$ws.Range("D41").Value = "'2.35"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.25%  "
$ws.Range("E42").Value = "  -2.34%  "
$ws.Range("E43").Value = "  -0.02%  "
$ws.Range("E44").Value = "  -4.14%  "
$ws.Range("E45").Value = "  +0.03%  "
$ws.Range("D46").Value = "'164.02"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.26%  "
$ws.Range("D47").Value = "'49.21"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.61%  "
$ws.Range("D48").Value = "'0.000301"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.62%  "
$ws.Range("D49").Value = "'418.89"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.92%  "
$ws.Range("D50").Value = "'1.38"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.05%  "
$ws.Range("D51").Value = "'8.61"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.99%  "
